$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("E3").Value = "  -5.42%  "
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("E6").Value = "  -7.79%  "
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("E8").Value = "  -11.00%  "
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("E14").Value = "  -6.53%  "
$ws.Range("E15").Value = "  -5.03%  "
$ws.Range("E16").Value = "  -6.35%  "
$ws.Range("E17").Value = "  -18.53%  "
$ws.Range("E18").Value = "  -9.36%  "
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("E22").Value = "  -6.91%  "
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("E25").Value = "  -10.27%  "
$ws.Range("E26").Value = "  +5.48%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("E29").Value = "  -4.90%  "
$ws.Range("E30").Value = "  -6.12%  "
$ws.Range("E31").Value = "  -8.29%  "
$ws.Range("E32").Value = "  -11.99%  "
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("E34").Value = "  -6.18%  "
$ws.Range("E35").Value = "  -7.92%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  -4.59%  "
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("E40").Value = "  -9.50%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("E47").Value = "  -10.44%  "
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("E49").Value = "  -11.09%  "
$ws.Range("E50").Value = "  -4.81%  "
